# Insert two new price-report rows right before the current row 83
# (Ají / Americana (o) - Primera and Segunda, dated 44546), shifting every
# subsequent data row down by two. This grows the sheet from A1:R194 to
# A1:R196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Insert()
$ws.Rows.Item(83).Insert()

# New row 83: Americana (o) / Primera
$ws.Range("A83").Value = 2
$ws.Range("B83").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44546
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 100112021
$ws.Range("G83").Value = "Ají"
$ws.Range("H83").Value = "Americana (o)"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 240
$ws.Range("K83").Value = 23000
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = 24000
$ws.Range("N83").Value = "$/caja 25 kilos"
$ws.Range("O83").Value = "Provincia de Limarí"
$ws.Range("P83").Value = 960
$ws.Range("Q83").Value = 25
$ws.Range("R83").Value = "Hortaliza"

# New row 84: Americana (o) / Segunda
$ws.Range("A84").Value = 2
$ws.Range("B84").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44546
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 100112021
$ws.Range("G84").Value = "Ají"
$ws.Range("H84").Value = "Americana (o)"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 13000
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = 14000
$ws.Range("N84").Value = "$/caja 25 kilos"
$ws.Range("O84").Value = "Provincia de Limarí"
$ws.Range("P84").Value = 560
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
